$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -22.06740000000001
$ws.Range("E3").Value = 16.2207

$ws.Range("A21").Value = -20.35069999999997

$ws.Range("A23").Value = -20.51959999999998
$ws.Range("E24").Value = 16.578

$ws.Range("A25").Value = -21.71819999999999

$ws.Range("B27").Value = 6.833000000000002

$ws.Range("B31").Value = 5.711800000000002

$ws.Range("B39").Value = 9.686000000000002

$ws.Range("B48").Value = 5.553900000000002

$ws.Range("B51").Value = 5.706699999999996

$ws.Range("B52").Value = 5.638699999999999

$ws.Range("A53").Value = -21.9205

$ws.Range("B55").Value = 6.129999999999995

$ws.Range("B56").Value = 5.277199999999997

$ws.Range("A57").Value = -22.63030000000002
$ws.Range("B57").Value = 4.840199999999995
$ws.Range("E57").Value = 16.59809999999999

$ws.Range("A59").Value = -22.1964

$ws.Range("E61").Value = 16.52930000000001

$ws.Range("A69").Value = -21.61519999999999

$ws.Range("E70").Value = 17.66630000000001

$ws.Range("B73").Value = 8.706899999999994

$ws.Range("A79").Value = -20.76890000000002

$ws.Range("A83").Value = -21.9547

$ws.Range("E86").Value = 16.71000000000002

$ws.Range("B89").Value = 4.982399999999995

$ws.Range("B90").Value = 5.935800000000004

$ws.Range("A93").Value = -21.3905

$ws.Range("E98").Value = 16.1046

$ws.Range("E100").Value = 16.41760000000001

$ws.Range("E102").Value = 16.61689999999998
